$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 (SamplesTab query) loses the Tumor / Analyte Type columns
$ws.Range("B3").Value = "SELECT`n    DISTINCT (smp.sample_id) AS `"Sample ID`",`n    sp.participant_id AS `"Participant ID`", `n    s.study_name AS `"Study Name`",`n    s.phs_accession AS Accession`nFROM `n    df_participant sp`nJOIN `n    df_study s ON sp.`"study.phs_accession`" = s.phs_accession`nJOIN `n    df_sample smp ON smp.`"participant.study_participant_id`" = sp.study_participant_id`nJOIN`n    df_diagnosis d ON d.`"participant.study_participant_id`" = sp.study_participant_id`nJOIN`n    df_program p ON p.program_acronym = s.`"program.program_acronym`"`nJOIN`n    df_file f1 ON f1.`"sample.sample_id`" = smp.sample_id`nJOIN`n    df_genomic_info gi ON gi.`"file.file_id`" = f1.file_id`nWHERE `n     s.phs_accession = 'phs002529' AND f1.experimental_strategy_and_data_subtypes = 'WGS'`nORDER BY `n    smp.sample_id ASC`nLIMIT 100;"

# Scroll / selection moved from C10 to B3
$ws.Range("B3").Select()
